$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header row - bold text on yellow fill
$ws.Range("A45").Value = "Ver0. problem"
$ws.Range("A45").Font.Bold = $true
$ws.Range("A45").Interior.Color = 65535

$ws.Range("A47").Value = "I2C pullup resistor position"
$ws.Range("C47").Value = "Position better close to MCU"

$ws.Range("A48").Value = "Mti-7 DK header position not correct"
$ws.Range("C48").Value = "1700mil->48.26mm"

$ws.Range("A49").Value = "Screw position blocked by LoRa"
$ws.Range("C49").Value = "put one LoRa in the middle of board"

$ws.Range("A50").Value = "Tempurature sensor ADC forgot"

$ws.Range("A46").Value = "I2C SDA SCL pullup resistor schematic false"
$ws.Range("C46").Value = "only need one set pullup resistor"

$ws.Range("A53").Value = "terminal, DCDC, through hole PAD not show up"

$ws.Range("A52").Value = "stm32, CAN controller,Top solder PAD not show up"
$ws.Range("B52").Value = "V"

$ws.Range("B53").Value = "V"

$ws.Range("B50").Select()
